$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = [double]"104.1546325"
$ws.Cells.Item(2, 8).Value = [double]"208.309265"
$ws.Cells.Item(2, 9).Value = [double]"0.9300139342750302"
$ws.Cells.Item(2, 10).Value = [double]"0.9008685911440201"
$ws.Cells.Item(2, 13).Value = [double]"62.009941"
$ws.Cells.Item(2, 14).Value = [double]"124.019882"
$ws.Cells.Item(2, 15).Value = [double]"0.9837206946349113"
$ws.Cells.Item(2, 16).Value = [double]"0.9775266241870173"
$ws.Cells.Item(2, 17).Value = [double]"6458.622616201682"
$ws.Cells.Item(2, 18).Value = [double]"25834.49046480673"
$ws.Cells.Item(2, 19).Value = [double]"0.9148739534451794"
$ws.Cells.Item(2, 20).Value = [double]"0.8806230327371283"

# Row 3
$ws.Cells.Item(3, 7).Value = [double]"104.1546325"
$ws.Cells.Item(3, 8).Value = [double]"208.309265"
$ws.Cells.Item(3, 9).Value = [double]"0.9300139342750302"
$ws.Cells.Item(3, 10).Value = [double]"0.9008685911440201"
$ws.Cells.Item(3, 15).Value = [double]"0.001352494296709523"
$ws.Cells.Item(3, 16).Value = [double]"0.002015967323812363"
$ws.Cells.Item(3, 17).Value = [double]"8.879807348419998"
$ws.Cells.Item(3, 18).Value = [double]"53.27884409052"
$ws.Cells.Item(3, 19).Value = [double]"0.001257838541967363"
$ws.Cells.Item(3, 20).Value = [double]"0.001816121642795224"

# Row 4
$ws.Cells.Item(4, 7).Value = [double]"104.1546325"
$ws.Cells.Item(4, 8).Value = [double]"208.309265"
$ws.Cells.Item(4, 9).Value = [double]"0.9300139342750302"
$ws.Cells.Item(4, 10).Value = [double]"0.9008685911440201"
$ws.Cells.Item(4, 11).Value = [double]"3"
$ws.Cells.Item(4, 12).Value = [double]"1"
$ws.Cells.Item(4, 13).Value = [double]"0.4639336666666667"
$ws.Cells.Item(4, 14).Value = [double]"1.391801"
$ws.Cells.Item(4, 15).Value = [double]"0.007359806209747155"
$ws.Cells.Item(4, 16).Value = [double]"0.01097019696462955"
$ws.Cells.Item(4, 17).Value = [double]"48.32084055604417"
$ws.Cells.Item(4, 18).Value = [double]"289.925043336265"
$ws.Cells.Item(4, 19).Value = [double]"0.00684472232862875"
$ws.Cells.Item(4, 20).Value = [double]"0.009882705884098229"

# Row 5
$ws.Cells.Item(5, 7).Value = [double]"104.1546325"
$ws.Cells.Item(5, 8).Value = [double]"208.309265"
$ws.Cells.Item(5, 9).Value = [double]"0.9300139342750302"
$ws.Cells.Item(5, 10).Value = [double]"0.9008685911440201"
$ws.Cells.Item(5, 13).Value = [double]"0.227331"
$ws.Cells.Item(5, 14).Value = [double]"0.454662"
$ws.Cells.Item(5, 15).Value = [double]"0.003606360619373094"
$ws.Cells.Item(5, 16).Value = [double]"0.003583652901767135"
$ws.Cells.Item(5, 17).Value = [double]"23.6775767608575"
$ws.Cells.Item(5, 18).Value = [double]"94.71030704342999"
$ws.Cells.Item(5, 19).Value = [double]"0.003353965628037706"
$ws.Cells.Item(5, 20).Value = [double]"0.003228400340764139"

# Row 6
$ws.Cells.Item(6, 7).Value = [double]"104.1546325"
$ws.Cells.Item(6, 8).Value = [double]"208.309265"
$ws.Cells.Item(6, 9).Value = [double]"0.9300139342750302"
$ws.Cells.Item(6, 10).Value = [double]"0.9008685911440201"
$ws.Cells.Item(6, 13).Value = [double]"0.1947683333333333"
$ws.Cells.Item(6, 14).Value = [double]"0.5843050000000001"
$ws.Cells.Item(6, 15).Value = [double]"0.003089789105904013"
$ws.Cells.Item(6, 16).Value = [double]"0.004605501028823711"
$ws.Cells.Item(6, 17).Value = [double]"20.28602418097083"
$ws.Cells.Item(6, 18).Value = [double]"121.716145085825"
$ws.Cells.Item(6, 19).Value = [double]"0.00287354692246192"
$ws.Cells.Item(6, 20).Value = [double]"0.004148951223348752"

# Row 7
$ws.Cells.Item(7, 7).Value = [double]"104.1546325"
$ws.Cells.Item(7, 8).Value = [double]"208.309265"
$ws.Cells.Item(7, 9).Value = [double]"0.9300139342750302"
$ws.Cells.Item(7, 10).Value = [double]"0.9008685911440201"
$ws.Cells.Item(7, 11).Value = [double]"2"
$ws.Cells.Item(7, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(7, 13).Value = [double]"0.05489533333333333"
$ws.Cells.Item(7, 14).Value = [double]"0.164686"
$ws.Cells.Item(7, 15).Value = [double]"0.0008708551333548545"
$ws.Cells.Item(7, 16).Value = [double]"0.001298057593949841"
$ws.Cells.Item(7, 17).Value = [double]"5.717603269298333"
$ws.Cells.Item(7, 18).Value = [double]"34.30561961579"
$ws.Cells.Item(7, 19).Value = [double]"0.0008099074087549544"
$ws.Cells.Item(7, 20).Value = [double]"0.001169379315885389"

# Row 8
$ws.Cells.Item(8, 9).Value = [double]"0.06205541323036482"
$ws.Cells.Item(8, 10).Value = [double]"0.09016602434017912"
$ws.Cells.Item(8, 13).Value = [double]"62.009941"
$ws.Cells.Item(8, 14).Value = [double]"124.019882"
$ws.Cells.Item(8, 15).Value = [double]"0.9837206946349113"
$ws.Cells.Item(8, 16).Value = [double]"0.9775266241870173"
$ws.Cells.Item(8, 17).Value = [double]"430.953215405104"
$ws.Cells.Item(8, 18).Value = [double]"2585.719292430624"
$ws.Cells.Item(8, 19).Value = [double]"0.06104519420883095"
$ws.Cells.Item(8, 20).Value = [double]"0.08813968938961973"

# Row 9
$ws.Cells.Item(9, 9).Value = [double]"0.06205541323036482"
$ws.Cells.Item(9, 10).Value = [double]"0.09016602434017912"
$ws.Cells.Item(9, 15).Value = [double]"0.001352494296709523"
$ws.Cells.Item(9, 16).Value = [double]"0.002015967323812363"
$ws.Cells.Item(9, 19).Value = [double]"8.392959247402109E-05"
$ws.Cells.Item(9, 20).Value = [double]"0.0001817717587878713"

# Row 10
$ws.Cells.Item(10, 9).Value = [double]"0.06205541323036482"
$ws.Cells.Item(10, 10).Value = [double]"0.09016602434017912"
$ws.Cells.Item(10, 11).Value = [double]"3"
$ws.Cells.Item(10, 12).Value = [double]"1"
$ws.Cells.Item(10, 13).Value = [double]"0.4639336666666667"
$ws.Cells.Item(10, 14).Value = [double]"1.391801"
$ws.Cells.Item(10, 15).Value = [double]"0.007359806209747155"
$ws.Cells.Item(10, 16).Value = [double]"0.01097019696462955"
$ws.Cells.Item(10, 17).Value = [double]"3.224220216314667"
$ws.Cells.Item(10, 18).Value = [double]"29.017981946832"
$ws.Cells.Item(10, 19).Value = [double]"0.0004567158156412648"
$ws.Cells.Item(10, 20).Value = [double]"0.0009891390465293472"

# Row 11
$ws.Cells.Item(11, 9).Value = [double]"0.06205541323036482"
$ws.Cells.Item(11, 10).Value = [double]"0.09016602434017912"
$ws.Cells.Item(11, 13).Value = [double]"0.227331"
$ws.Cells.Item(11, 14).Value = [double]"0.454662"
$ws.Cells.Item(11, 15).Value = [double]"0.003606360619373094"
$ws.Cells.Item(11, 16).Value = [double]"0.003583652901767135"
$ws.Cells.Item(11, 17).Value = [double]"1.579892253264"
$ws.Cells.Item(11, 18).Value = [double]"9.479353519584"
$ws.Cells.Item(11, 19).Value = [double]"0.0002237941984929118"
$ws.Cells.Item(11, 20).Value = [double]"0.0003231237347674891"

# Row 12
$ws.Cells.Item(12, 9).Value = [double]"0.06205541323036482"
$ws.Cells.Item(12, 10).Value = [double]"0.09016602434017912"
$ws.Cells.Item(12, 13).Value = [double]"0.1947683333333333"
$ws.Cells.Item(12, 14).Value = [double]"0.5843050000000001"
$ws.Cells.Item(12, 15).Value = [double]"0.003089789105904013"
$ws.Cells.Item(12, 16).Value = [double]"0.004605501028823711"
$ws.Cells.Item(12, 17).Value = [double]"1.353590055973333"
$ws.Cells.Item(12, 18).Value = [double]"12.18231050376"
$ws.Cells.Item(12, 19).Value = [double]"0.000191738139761553"
$ws.Cells.Item(12, 20).Value = [double]"0.0004152597178636387"

# Row 13
$ws.Cells.Item(13, 9).Value = [double]"0.06205541323036482"
$ws.Cells.Item(13, 10).Value = [double]"0.09016602434017912"
$ws.Cells.Item(13, 11).Value = [double]"2"
$ws.Cells.Item(13, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(13, 13).Value = [double]"0.05489533333333333"
$ws.Cells.Item(13, 14).Value = [double]"0.164686"
$ws.Cells.Item(13, 15).Value = [double]"0.0008708551333548545"
$ws.Cells.Item(13, 16).Value = [double]"0.001298057593949841"
$ws.Cells.Item(13, 17).Value = [double]"0.3815085134613333"
$ws.Cells.Item(13, 18).Value = [double]"3.433576621152"
$ws.Cells.Item(13, 19).Value = [double]"5.404127516411996E-05"
$ws.Cells.Item(13, 20).Value = [double]"0.0001170406926110357"

# Row 14
$ws.Cells.Item(14, 7).Value = [double]"0.591442"
$ws.Cells.Item(14, 8).Value = [double]"1.182884"
$ws.Cells.Item(14, 9).Value = [double]"0.005281083405632414"
$ws.Cells.Item(14, 10).Value = [double]"0.00511558159723142"
$ws.Cells.Item(14, 13).Value = [double]"62.009941"
$ws.Cells.Item(14, 14).Value = [double]"124.019882"
$ws.Cells.Item(14, 15).Value = [double]"0.9837206946349113"
$ws.Cells.Item(14, 16).Value = [double]"0.9775266241870173"
$ws.Cells.Item(14, 17).Value = [double]"36.675283524922"
$ws.Cells.Item(14, 18).Value = [double]"146.701134099688"
$ws.Cells.Item(14, 19).Value = [double]"0.005195111036213621"
$ws.Cells.Item(14, 20).Value = [double]"0.00500061720949486"

# Row 15
$ws.Cells.Item(15, 7).Value = [double]"0.591442"
$ws.Cells.Item(15, 8).Value = [double]"1.182884"
$ws.Cells.Item(15, 9).Value = [double]"0.005281083405632414"
$ws.Cells.Item(15, 10).Value = [double]"0.00511558159723142"
$ws.Cells.Item(15, 15).Value = [double]"0.001352494296709523"
$ws.Cells.Item(15, 16).Value = [double]"0.002015967323812363"
$ws.Cells.Item(15, 17).Value = [double]"0.050423979152"
$ws.Cells.Item(15, 18).Value = [double]"0.302543874912"
$ws.Cells.Item(15, 19).Value = [double]"7.142635186565142E-06"
$ws.Cells.Item(15, 20).Value = [double]"1.03128453423144E-05"

# Row 16
$ws.Cells.Item(16, 7).Value = [double]"0.591442"
$ws.Cells.Item(16, 8).Value = [double]"1.182884"
$ws.Cells.Item(16, 9).Value = [double]"0.005281083405632414"
$ws.Cells.Item(16, 10).Value = [double]"0.00511558159723142"
$ws.Cells.Item(16, 11).Value = [double]"3"
$ws.Cells.Item(16, 12).Value = [double]"1"
$ws.Cells.Item(16, 13).Value = [double]"0.4639336666666667"
$ws.Cells.Item(16, 14).Value = [double]"1.391801"
$ws.Cells.Item(16, 15).Value = [double]"0.007359806209747155"
$ws.Cells.Item(16, 16).Value = [double]"0.01097019696462955"
$ws.Cells.Item(16, 17).Value = [double]"0.2743898556806667"
$ws.Cells.Item(16, 18).Value = [double]"1.646339134084"
$ws.Cells.Item(16, 19).Value = [double]"3.88677504429661E-05"
$ws.Cells.Item(16, 20).Value = [double]"5.611893771026291E-05"

# Row 17
$ws.Cells.Item(17, 7).Value = [double]"0.591442"
$ws.Cells.Item(17, 8).Value = [double]"1.182884"
$ws.Cells.Item(17, 9).Value = [double]"0.005281083405632414"
$ws.Cells.Item(17, 10).Value = [double]"0.00511558159723142"
$ws.Cells.Item(17, 13).Value = [double]"0.227331"
$ws.Cells.Item(17, 14).Value = [double]"0.454662"
$ws.Cells.Item(17, 15).Value = [double]"0.003606360619373094"
$ws.Cells.Item(17, 16).Value = [double]"0.003583652901767135"
$ws.Cells.Item(17, 17).Value = [double]"0.134453101302"
$ws.Cells.Item(17, 18).Value = [double]"0.537812405208"
$ws.Cells.Item(17, 19).Value = [double]"1.904549122169748E-05"
$ws.Cells.Item(17, 20).Value = [double]"1.833246883514493E-05"

# Row 18
$ws.Cells.Item(18, 7).Value = [double]"0.591442"
$ws.Cells.Item(18, 8).Value = [double]"1.182884"
$ws.Cells.Item(18, 9).Value = [double]"0.005281083405632414"
$ws.Cells.Item(18, 10).Value = [double]"0.00511558159723142"
$ws.Cells.Item(18, 13).Value = [double]"0.1947683333333333"
$ws.Cells.Item(18, 14).Value = [double]"0.5843050000000001"
$ws.Cells.Item(18, 15).Value = [double]"0.003089789105904013"
$ws.Cells.Item(18, 16).Value = [double]"0.004605501028823711"
$ws.Cells.Item(18, 17).Value = [double]"0.1151941726033334"
$ws.Cells.Item(18, 18).Value = [double]"0.6911650356200001"
$ws.Cells.Item(18, 19).Value = [double]"1.63174339740935E-05"
$ws.Cells.Item(18, 20).Value = [double]"2.355981630908095E-05"

# Row 19
$ws.Cells.Item(19, 7).Value = [double]"0.591442"
$ws.Cells.Item(19, 8).Value = [double]"1.182884"
$ws.Cells.Item(19, 9).Value = [double]"0.005281083405632414"
$ws.Cells.Item(19, 10).Value = [double]"0.00511558159723142"
$ws.Cells.Item(19, 11).Value = [double]"2"
$ws.Cells.Item(19, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(19, 13).Value = [double]"0.05489533333333333"
$ws.Cells.Item(19, 14).Value = [double]"0.164686"
$ws.Cells.Item(19, 15).Value = [double]"0.0008708551333548545"
$ws.Cells.Item(19, 16).Value = [double]"0.001298057593949841"
$ws.Cells.Item(19, 17).Value = [double]"0.03246740573733333"
$ws.Cells.Item(19, 18).Value = [double]"0.194804434424"
$ws.Cells.Item(19, 19).Value = [double]"4.599058593470125E-06"
$ws.Cells.Item(19, 20).Value = [double]"6.640319539756299E-06"

# Row 20
$ws.Cells.Item(20, 5).Value = [double]"1"
$ws.Cells.Item(20, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(20, 7).Value = [double]"0.296732"
$ws.Cells.Item(20, 8).Value = [double]"0.890196"
$ws.Cells.Item(20, 9).Value = [double]"0.002649569088972574"
$ws.Cells.Item(20, 10).Value = [double]"0.003849802918569379"
$ws.Cells.Item(20, 13).Value = [double]"62.009941"
$ws.Cells.Item(20, 14).Value = [double]"124.019882"
$ws.Cells.Item(20, 15).Value = [double]"0.9837206946349113"
$ws.Cells.Item(20, 16).Value = [double]"0.9775266241870173"
$ws.Cells.Item(20, 17).Value = [double]"18.400333812812"
$ws.Cells.Item(20, 18).Value = [double]"110.402002876872"
$ws.Cells.Item(20, 19).Value = [double]"0.00260643594468729"
$ws.Cells.Item(20, 20).Value = [double]"0.003763284850774451"

# Row 21
$ws.Cells.Item(21, 5).Value = [double]"1"
$ws.Cells.Item(21, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(21, 7).Value = [double]"0.296732"
$ws.Cells.Item(21, 8).Value = [double]"0.890196"
$ws.Cells.Item(21, 9).Value = [double]"0.002649569088972574"
$ws.Cells.Item(21, 10).Value = [double]"0.003849802918569379"
$ws.Cells.Item(21, 15).Value = [double]"0.001352494296709523"
$ws.Cells.Item(21, 16).Value = [double]"0.002015967323812363"
$ws.Cells.Item(21, 17).Value = [double]"0.025298183392"
$ws.Cells.Item(21, 18).Value = [double]"0.227683650528"
$ws.Cells.Item(21, 19).Value = [double]"3.583527081573252E-06"
$ws.Cells.Item(21, 20).Value = [double]"7.761076886953336E-06"

# Row 22
$ws.Cells.Item(22, 5).Value = [double]"1"
$ws.Cells.Item(22, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(22, 7).Value = [double]"0.296732"
$ws.Cells.Item(22, 8).Value = [double]"0.890196"
$ws.Cells.Item(22, 9).Value = [double]"0.002649569088972574"
$ws.Cells.Item(22, 10).Value = [double]"0.003849802918569379"
$ws.Cells.Item(22, 11).Value = [double]"3"
$ws.Cells.Item(22, 12).Value = [double]"1"
$ws.Cells.Item(22, 13).Value = [double]"0.4639336666666667"
$ws.Cells.Item(22, 14).Value = [double]"1.391801"
$ws.Cells.Item(22, 15).Value = [double]"0.007359806209747155"
$ws.Cells.Item(22, 16).Value = [double]"0.01097019696462955"
$ws.Cells.Item(22, 17).Value = [double]"0.1376639647773333"
$ws.Cells.Item(22, 18).Value = [double]"1.238975682996"
$ws.Cells.Item(22, 19).Value = [double]"1.950031503417446E-05"
$ws.Cells.Item(22, 20).Value = [double]"4.223309629171179E-05"

# Row 23
$ws.Cells.Item(23, 5).Value = [double]"1"
$ws.Cells.Item(23, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(23, 7).Value = [double]"0.296732"
$ws.Cells.Item(23, 8).Value = [double]"0.890196"
$ws.Cells.Item(23, 9).Value = [double]"0.002649569088972574"
$ws.Cells.Item(23, 10).Value = [double]"0.003849802918569379"
$ws.Cells.Item(23, 13).Value = [double]"0.227331"
$ws.Cells.Item(23, 14).Value = [double]"0.454662"
$ws.Cells.Item(23, 15).Value = [double]"0.003606360619373094"
$ws.Cells.Item(23, 16).Value = [double]"0.003583652901767135"
$ws.Cells.Item(23, 17).Value = [double]"0.067456382292"
$ws.Cells.Item(23, 18).Value = [double]"0.404738293752"
$ws.Cells.Item(23, 19).Value = [double]"9.555301620778937E-06"
$ws.Cells.Item(23, 20).Value = [double]"1.379635740036274E-05"

# Row 24
$ws.Cells.Item(24, 5).Value = [double]"1"
$ws.Cells.Item(24, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(24, 7).Value = [double]"0.296732"
$ws.Cells.Item(24, 8).Value = [double]"0.890196"
$ws.Cells.Item(24, 9).Value = [double]"0.002649569088972574"
$ws.Cells.Item(24, 10).Value = [double]"0.003849802918569379"
$ws.Cells.Item(24, 13).Value = [double]"0.1947683333333333"
$ws.Cells.Item(24, 14).Value = [double]"0.5843050000000001"
$ws.Cells.Item(24, 15).Value = [double]"0.003089789105904013"
$ws.Cells.Item(24, 16).Value = [double]"0.004605501028823711"
$ws.Cells.Item(24, 17).Value = [double]"0.05779399708666667"
$ws.Cells.Item(24, 18).Value = [double]"0.52014597378"
$ws.Cells.Item(24, 19).Value = [double]"8.186609706447481E-06"
$ws.Cells.Item(24, 20).Value = [double]"1.77302713022398E-05"

# Row 25
$ws.Cells.Item(25, 5).Value = [double]"1"
$ws.Cells.Item(25, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(25, 7).Value = [double]"0.296732"
$ws.Cells.Item(25, 8).Value = [double]"0.890196"
$ws.Cells.Item(25, 9).Value = [double]"0.002649569088972574"
$ws.Cells.Item(25, 10).Value = [double]"0.003849802918569379"
$ws.Cells.Item(25, 11).Value = [double]"2"
$ws.Cells.Item(25, 12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(25, 13).Value = [double]"0.05489533333333333"
$ws.Cells.Item(25, 14).Value = [double]"0.164686"
$ws.Cells.Item(25, 15).Value = [double]"0.0008708551333548545"
$ws.Cells.Item(25, 16).Value = [double]"0.001298057593949841"
$ws.Cells.Item(25, 17).Value = [double]"0.01628920205066666"
$ws.Cells.Item(25, 18).Value = [double]"0.146602818456"
$ws.Cells.Item(25, 19).Value = [double]"2.307390842310112E-06"
$ws.Cells.Item(25, 20).Value = [double]"4.997265913659241E-06"

